# Add a new paragraph after "URL to GitHub Repository:" containing the
# repository URL, matching the bold formatting of the preceding paragraph.

$d = $word.ActiveDocument

# Locate the existing "URL to GitHub Repository:" paragraph and grab the
# end of its range so we can append a new paragraph right after it.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "URL to GitHub Repository:") {
        $target = $p
    }
}

$endRange = $target.Range
$endRange.Collapse(0)  # wdCollapseEnd
$endRange.InsertParagraphAfter()

# Move into the freshly-created paragraph and set its text/formatting.
$newRange = $target.Next().Range
$newRange.Text = "https://github.com/wensun163/backendbootcampweek10"
$newRange.Font.Bold = $true
